# Cminox productos.xlsx - add image URLs for "Unión H/H" and "Tapón Cachucha" rows
# Commit: "Se agrega imagen de Reducción Campana"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$unionUrl    = "https://raw.githubusercontent.com/Rick130425/CminoxImages/main/Uni%C3%B3n%20H-H/union-hh.jpg,http://http2.mlstatic.com/D_743692-MLM47153394760_082021-O.jpg"
$cachuchaUrl = "https://raw.githubusercontent.com/Rick130425/CminoxImages/main/Tap%C3%B3n%20Cachuca/Tap%C3%B3n%20Cachucha.jpg,http://http2.mlstatic.com/D_743692-MLM47153394760_082021-O.jpg"

# Rows 2-10 (Unión H/H 1/4 .. 2 1/2) already show an image link - just refresh
# the displayed text to the new "union-hh" image (keeps the existing hyperlink).
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("G$r").Value = $unionUrl
}

# Rows 11-12 (Unión H/H 3, Unión H/H 4) had no image before - add one now.
foreach ($r in 11, 12) {
    $cell = $ws.Range("G$r")
    $ws.Hyperlinks.Add($cell, $unionUrl, "", "", $unionUrl) | Out-Null
    $cell.Style = $ws.Range("G2").Style
}

# Rows 13-23 (Tapón Cachucha 1/4 .. 4) had no image before - add the cap photo.
for ($r = 13; $r -le 23; $r++) {
    $cell = $ws.Range("G$r")
    $ws.Hyperlinks.Add($cell, $cachuchaUrl, "", "", $cachuchaUrl) | Out-Null
    $cell.Style = $ws.Range("G2").Style
}

# Move the active selection from G13 to A13, as in the saved workbook.
$ws.Range("A13").Select()
